# Update the "想去人数" (interested-count) figures to the latest scraped
# values, as produced by the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" holds the exhibition listing; row 3 is the 2024 良苗动漫
# autumn gala, row 10 is the Lolita tea party.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5989
$wsExhibit.Range("F10").Value = 29

# Sheet "全部类型" aggregates every event type; the same two events live
# at row 3 and row 12 there.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5989
$wsAll.Range("F12").Value = 29
